$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data row contents (row 2), keeping formatting on A2
$ws.Range("A2:C2").ClearContents()

# Match the saved selection state (active cell H8)
[void]$ws.Range("H8").Select()
